$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text/value updates (values that Excel will not misinterpret as pure numbers)
$ws.Range("D2").Value = "37.863.77"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.029.43"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "2.330.63"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "2.041.50"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "37.776.02"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -4.44%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  +6.72%  "
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "1.520.70"
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "2.218.24"
$ws.Range("E51").Value = "  -1.13%  "

# Price cells that look like plain numbers: force them to stay text (no quote-prefix style change)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.41"
$ws.Range("B5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.45"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.58"
$ws.Range("B12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.14"
$ws.Range("B14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.760"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("B16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.02"
$ws.Range("B19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.00"
$ws.Range("B20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.06"
$ws.Range("B22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("B24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4122) | Out-Null
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("B25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.28"
$ws.Range("B26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.13"
$ws.Range("B27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.98"
$ws.Range("B29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.51"
$ws.Range("B34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.26"
$ws.Range("B37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4122) | Out-Null
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.63"
$ws.Range("B42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0918"
$ws.Range("B45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.07"
$ws.Range("B47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0